$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A keeps its original custom width (15.42578125); re-assert it so the
# engine does not clear its customWidth flag when column B's width changes.
# Column B's width changes from 14.7109375 to 15.42578125 (character width
# units round to the nearest 1/6 in this engine, so 14.6667 -> 15.5 is the
# closest achievable approximation of 15.42578125).
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666

$ws.Range("A1").Value = -0.21729837854633161
$ws.Range("B1").Value = 0.21715215758334239
$ws.Range("A2").Value = -0.13853351630225053
$ws.Range("B2").Value = 0.13821373841016005
$ws.Range("A3").Value = -0.088513061909686996
$ws.Range("B3").Value = 0.088198351843887579
$ws.Range("A4").Value = -0.080198351906064502
$ws.Range("B4").Value = 0.079867768381042481
$ws.Range("A5").Value = -0.076867768415365134
$ws.Range("B5").Value = 0.075756255420880159
$ws.Range("A6").Value = -0.0035608320954096939
$ws.Range("B6").Value = 0.0034955014142390439
$ws.Range("A7").Value = 0.0065044984977862264
$ws.Range("B7").Value = -0.0065111306408738301
$ws.Range("A8").Value = 0.01651113055318465
$ws.Range("B8").Value = -0.016523983919472851
$ws.Range("A9").Value = -0.038715522649068745
$ws.Range("B9").Value = 0.0385753608146997
$ws.Range("A10").Value = -0.036575360856597072
$ws.Range("B10").Value = 0.036566926066562644
$ws.Range("A11").Value = -0.033566926115027762
$ws.Range("B11").Value = 0.03355254834084942
$ws.Range("A12").Value = -0.030052548393217027
$ws.Range("B12").Value = 0.029950536357071122
$ws.Range("A13").Value = -0.017164167583041845
$ws.Range("B13").Value = 0.017079089337326714
$ws.Range("A14").Value = -0.0090790894243033549
$ws.Range("B14").Value = 0.009051366959320184
$ws.Range("A15").Value = -0.0080513670035058382
$ws.Range("B15").Value = 0.0080334642346979734
$ws.Range("A16").Value = -0.006033464285689405
$ws.Range("B16").Value = 0.0060033594696080606
$ws.Range("A17").Value = -0.0040033595215875906
$ws.Range("B17").Value = 0.0039999999353890203
$ws.Range("A18").Value = -0.072618641319497357
$ws.Range("B18").Value = 0.072492798748438503
$ws.Range("A19").Value = -0.068492798775567909
$ws.Range("B19").Value = 0.067550921736479363
$ws.Range("A20").Value = -0.008016965842722712
$ws.Range("B20").Value = 0.0080057589827813302
$ws.Range("A21").Value = -0.004005759024662936
$ws.Range("B21").Value = 0.0039999999578590462
$ws.Range("A22").Value = -0.045700676534343287
$ws.Range("B22").Value = 0.045490824603374946
$ws.Range("A23").Value = -0.040490824645186052
$ws.Range("B23").Value = 0.040097358410300821
$ws.Range("A24").Value = -0.020097358549171496
$ws.Range("B24").Value = 0.019999999859209971
$ws.Range("A25").Value = -0.070195423363598408
$ws.Range("B25").Value = 0.070110626207229743
$ws.Range("A26").Value = -0.06761062624994274
$ws.Range("B26").Value = 0.067503296026490744
$ws.Range("A27").Value = -0.065003296071147521
$ws.Range("B27").Value = 0.064379930344526493
$ws.Range("A28").Value = -0.062379930395991323
$ws.Range("B28").Value = 0.061965362732957274
$ws.Range("A29").Value = -0.054965362822138708
$ws.Range("B29").Value = 0.054854100186421562
$ws.Range("A30").Value = 0.0051458993925415086
$ws.Range("B30").Value = -0.0052439431615844967
$ws.Range("A31").Value = 0.012243943069389474
$ws.Range("B31").Value = -0.01225986554154268
$ws.Range("A32").Value = -0.0040009905243216792
$ws.Range("B32").Value = 0.0039999999265916131
